$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4.5883315914771696
$ws.Range("C2").Value = 0.021196651826595199
$ws.Range("D2").Value = 0.80794778204261997

[void]$ws.Range("D2").Select()
